$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 11.22947716699574
$ws.Cells.Item(2, 3).Value = 4.562841544432813
$ws.Cells.Item(2, 5).Value = 21.262535195682
$ws.Cells.Item(2, 6).Value = 42.88119202077788
$ws.Cells.Item(2, 7).Value = 35.55100753008026
$ws.Cells.Item(2, 8).Value = 15.96393164064879
$ws.Cells.Item(2, 10).Value = 8.496830299205076
$ws.Cells.Item(2, 11).Value = 10.57308380515595
$ws.Cells.Item(2, 14).Value = 19.18376405186429

$ws.Cells.Item(3, 2).Value = 10.96837659746983
$ws.Cells.Item(3, 3).Value = 4.322609844355632
$ws.Cells.Item(3, 5).Value = 20.97705558039033
$ws.Cells.Item(3, 6).Value = 42.67627314368968
$ws.Cells.Item(3, 7).Value = 35.59053336244161
$ws.Cells.Item(3, 8).Value = 16.01560357352362
$ws.Cells.Item(3, 10).Value = 8.521570939434906
$ws.Cells.Item(3, 11).Value = 10.39788207326725
$ws.Cells.Item(3, 14).Value = 19.25318451558792

$ws.Cells.Item(4, 2).Value = 10.80761167183282
$ws.Cells.Item(4, 3).Value = 4.167219976916611
$ws.Cells.Item(4, 5).Value = 20.80503561848765
$ws.Cells.Item(4, 6).Value = 42.56207973144527
$ws.Cells.Item(4, 7).Value = 35.6266108291538
$ws.Cells.Item(4, 8).Value = 16.050222720301
$ws.Cells.Item(4, 10).Value = 8.537764533557279
$ws.Cells.Item(4, 11).Value = 10.29118706943581
$ws.Cells.Item(4, 14).Value = 19.29772476309124

$ws.Cells.Item(5, 2).Value = 10.74208870475931
$ws.Cells.Item(5, 3).Value = 4.10194538026719
$ws.Cells.Item(5, 5).Value = 20.73583965195489
$ws.Cells.Item(5, 6).Value = 42.51850253178578
$ws.Cells.Item(5, 7).Value = 35.6442687742647
$ws.Cells.Item(5, 8).Value = 16.0650564387273
$ws.Cells.Item(5, 10).Value = 8.544616017188014
$ws.Cells.Item(5, 11).Value = 10.24798702461825
$ws.Cells.Item(5, 14).Value = 19.31635858341925

$ws.Cells.Item(6, 2).Value = 10.73121112904774
$ws.Cells.Item(6, 3).Value = 4.0909897337304
$ws.Cells.Item(6, 5).Value = 20.7244066603536
$ws.Cells.Item(6, 6).Value = 42.51144609323688
$ws.Cells.Item(6, 7).Value = 35.6473789795797
$ws.Cells.Item(6, 8).Value = 16.06756338859538
$ws.Cells.Item(6, 10).Value = 8.545768958320313
$ws.Cells.Item(6, 11).Value = 10.24083223495622
$ws.Cells.Item(6, 14).Value = 19.31948195163963

$ws.Cells.Item(7, 2).Value = 10.80672790841547
$ws.Cells.Item(7, 3).Value = 4.166347518182844
$ws.Cells.Item(7, 5).Value = 20.80409865145973
$ws.Cells.Item(7, 6).Value = 42.56148001815146
$ws.Cells.Item(7, 7).Value = 35.62683702028181
$ws.Cells.Item(7, 8).Value = 16.05041983451593
$ws.Cells.Item(7, 10).Value = 8.537855912428927
$ws.Cells.Item(7, 11).Value = 10.29060325224972
$ws.Cells.Item(7, 14).Value = 19.29797410636416

$ws.Cells.Item(8, 2).Value = 11.13960327838751
$ws.Cells.Item(8, 3).Value = 4.481663246926431
$ws.Cells.Item(8, 5).Value = 21.16347394723262
$ws.Cells.Item(8, 6).Value = 42.80814427451651
$ws.Cells.Item(8, 7).Value = 35.56217778356081
$ws.Cells.Item(8, 8).Value = 15.98114710882555
$ws.Cells.Item(8, 10).Value = 8.505152942042765
$ws.Cells.Item(8, 11).Value = 10.51252635343783
$ws.Cells.Item(8, 14).Value = 19.2073035494454

$ws.Cells.Item(9, 2).Value = 11.78448862625567
$ws.Cells.Item(9, 3).Value = 5.036358224527242
$ws.Cells.Item(9, 5).Value = 21.89043410500549
$ws.Cells.Item(9, 6).Value = 43.38242029389529
$ws.Cells.Item(9, 7).Value = 35.52958299746476
$ws.Cells.Item(9, 8).Value = 15.86830554892983
$ws.Cells.Item(9, 10).Value = 8.448964652320399
$ws.Cells.Item(9, 11).Value = 10.95217131393157
$ws.Cells.Item(9, 14).Value = 19.04462715188334

$ws.Cells.Item(10, 2).Value = 12.24800264078519
$ws.Cells.Item(10, 3).Value = 5.404031519395515
$ws.Cells.Item(10, 5).Value = 22.43310490719471
$ws.Cells.Item(10, 6).Value = 43.85699765447095
$ws.Cells.Item(10, 7).Value = 35.56362911464642
$ws.Cells.Item(10, 8).Value = 15.79949492768069
$ws.Cells.Item(10, 10).Value = 8.412505353635343
$ws.Cells.Item(10, 11).Value = 11.27457271568616
$ws.Cells.Item(10, 14).Value = 18.93422761889357

$ws.Cells.Item(11, 2).Value = 12.45557844771265
$ws.Cells.Item(11, 3).Value = 5.562495889572937
$ws.Cells.Item(11, 5).Value = 22.6807908816247
$ws.Cells.Item(11, 6).Value = 44.08370166143813
$ws.Cells.Item(11, 7).Value = 35.59178739068661
$ws.Cells.Item(11, 8).Value = 15.77126674213665
$ws.Cells.Item(11, 10).Value = 8.3969621751798
$ws.Cells.Item(11, 11).Value = 11.4204337569454
$ws.Cells.Item(11, 14).Value = 18.8859618473296

$ws.Cells.Item(12, 2).Value = 12.53363023181476
$ws.Cells.Item(12, 3).Value = 5.621229661519625
$ws.Cells.Item(12, 5).Value = 22.77461518894668
$ws.Cells.Item(12, 6).Value = 44.17104570828995
$ws.Cells.Item(12, 7).Value = 35.6042748372781
$ws.Cells.Item(12, 8).Value = 15.76102077116995
$ws.Cells.Item(12, 10).Value = 8.391225985607495
$ws.Cells.Item(12, 11).Value = 11.47549854712056
$ws.Cells.Item(12, 14).Value = 18.86796440845463

$ws.Cells.Item(13, 2).Value = 12.51684619258652
$ws.Cells.Item(13, 3).Value = 5.608637037318466
$ws.Cells.Item(13, 5).Value = 22.75440843404624
$ws.Cells.Item(13, 6).Value = 44.15216903855445
$ws.Cells.Item(13, 7).Value = 35.60150428066969
$ws.Cells.Item(13, 8).Value = 15.76320768093192
$ws.Cells.Item(13, 10).Value = 8.392454723129882
$ws.Cells.Item(13, 11).Value = 11.4636477452781
$ws.Cells.Item(13, 14).Value = 18.87182805931528

$ws.Cells.Item(14, 2).Value = 12.46201133496878
$ws.Cells.Item(14, 3).Value = 5.567353512615743
$ws.Cells.Item(14, 5).Value = 22.68850973497625
$ws.Cells.Item(14, 6).Value = 44.09085787366773
$ws.Cells.Item(14, 7).Value = 35.59277816813474
$ws.Cells.Item(14, 8).Value = 15.77041490511878
$ws.Cells.Item(14, 10).Value = 8.396487256891943
$ws.Cells.Item(14, 11).Value = 11.42496767339337
$ws.Cells.Item(14, 14).Value = 18.88447558859183

$ws.Cells.Item(15, 2).Value = 12.4283491205249
$ws.Cells.Item(15, 3).Value = 5.541900124074518
$ws.Cells.Item(15, 5).Value = 22.64814636176577
$ws.Cells.Item(15, 6).Value = 44.05349604837345
$ws.Cells.Item(15, 7).Value = 35.58767081382379
$ws.Cells.Item(15, 8).Value = 15.77488732628448
$ws.Cells.Item(15, 10).Value = 8.398976786175846
$ws.Cells.Item(15, 11).Value = 11.40125134666777
$ws.Cells.Item(15, 14).Value = 18.89225895602722

$ws.Cells.Item(16, 2).Value = 12.23436401897217
$ws.Cells.Item(16, 3).Value = 5.393497757593669
$ws.Cells.Item(16, 5).Value = 22.41692709111976
$ws.Cells.Item(16, 6).Value = 43.84239467722432
$ws.Cells.Item(16, 7).Value = 35.56204413551568
$ws.Cells.Item(16, 8).Value = 15.80140168308372
$ws.Cells.Item(16, 10).Value = 8.413542093086694
$ws.Cells.Item(16, 11).Value = 11.2650195300109
$ws.Cells.Item(16, 14).Value = 18.93742112738454

$ws.Cells.Item(17, 2).Value = 12.11446098693924
$ws.Cells.Item(17, 3).Value = 5.300198659474701
$ws.Cells.Item(17, 5).Value = 22.27522949932229
$ws.Cells.Item(17, 6).Value = 43.71562003745429
$ws.Cells.Item(17, 7).Value = 35.54957073786782
$ws.Cells.Item(17, 8).Value = 15.81845572127088
$ws.Cells.Item(17, 10).Value = 8.422744251423836
$ws.Cells.Item(17, 11).Value = 11.18120112296359
$ws.Cells.Item(17, 14).Value = 18.96562648605178

$ws.Cells.Item(18, 2).Value = 12.04519177441898
$ws.Cells.Item(18, 3).Value = 5.24570829297206
$ws.Cells.Item(18, 5).Value = 22.19381130558628
$ws.Cells.Item(18, 6).Value = 43.64372446416675
$ws.Cells.Item(18, 7).Value = 35.54358897334017
$ws.Cells.Item(18, 8).Value = 15.82855402152855
$ws.Cells.Item(18, 10).Value = 8.428135211577992
$ws.Cells.Item(18, 11).Value = 11.13291872060053
$ws.Cells.Item(18, 14).Value = 18.98203361151765

$ws.Cells.Item(19, 2).Value = 12.02168873910267
$ws.Cells.Item(19, 3).Value = 5.227117034955365
$ws.Cells.Item(19, 5).Value = 22.16626136544727
$ws.Cells.Item(19, 6).Value = 43.61955911405613
$ws.Cells.Item(19, 7).Value = 35.54176836760097
$ws.Cells.Item(19, 8).Value = 15.8320227724244
$ws.Cells.Item(19, 10).Value = 8.429977357230364
$ws.Cells.Item(19, 11).Value = 11.1165604121013
$ws.Cells.Item(19, 14).Value = 18.98762045419999

$ws.Cells.Item(20, 2).Value = 12.12725695539934
$ws.Cells.Item(20, 3).Value = 5.310216191976208
$ws.Cells.Item(20, 5).Value = 22.29030554823512
$ws.Cells.Item(20, 6).Value = 43.72901007594027
$ws.Cells.Item(20, 7).Value = 35.55077509674965
$ws.Cells.Item(20, 8).Value = 15.81661034151721
$ws.Cells.Item(20, 10).Value = 8.42175451225717
$ws.Cells.Item(20, 11).Value = 11.19013161694332
$ws.Cells.Item(20, 14).Value = 18.96260492916864

$ws.Cells.Item(21, 2).Value = 12.47813327840403
$ws.Cells.Item(21, 3).Value = 5.579514099218754
$ws.Cells.Item(21, 5).Value = 22.70786562406019
$ws.Cells.Item(21, 6).Value = 44.10882633092675
$ws.Cells.Item(21, 7).Value = 35.59529171335011
$ws.Cells.Item(21, 8).Value = 15.76828592321843
$ws.Cells.Item(21, 10).Value = 8.395298743495401
$ws.Cells.Item(21, 11).Value = 11.43633397392154
$ws.Cells.Item(21, 14).Value = 18.88075312145884

$ws.Cells.Item(22, 2).Value = 12.70419641479721
$ws.Cells.Item(22, 3).Value = 5.748091901632868
$ws.Cells.Item(22, 5).Value = 22.98090705261716
$ws.Cells.Item(22, 6).Value = 44.36574940478897
$ws.Cells.Item(22, 7).Value = 35.63502012712112
$ws.Cells.Item(22, 8).Value = 15.73928839835076
$ws.Cells.Item(22, 10).Value = 8.378880695943248
$ws.Cells.Item(22, 11).Value = 11.59622919884312
$ws.Cells.Item(22, 14).Value = 18.82888823450064

$ws.Cells.Item(23, 2).Value = 12.58386535428238
$ws.Cells.Item(23, 3).Value = 5.658800386730934
$ws.Cells.Item(23, 5).Value = 22.83519500600293
$ws.Cells.Item(23, 6).Value = 44.22784982456472
$ws.Cells.Item(23, 7).Value = 35.61284305620814
$ws.Cells.Item(23, 8).Value = 15.75452792031342
$ws.Cells.Item(23, 10).Value = 8.387563566902804
$ws.Cells.Item(23, 11).Value = 11.51100010438656
$ws.Cells.Item(23, 14).Value = 18.856420829331

$ws.Cells.Item(24, 2).Value = 12.12147293637296
$ws.Cells.Item(24, 3).Value = 5.305689917038553
$ws.Cells.Item(24, 5).Value = 22.28348951768141
$ws.Cells.Item(24, 6).Value = 43.72295335471283
$ws.Cells.Item(24, 7).Value = 35.55022690140844
$ws.Cells.Item(24, 8).Value = 15.81744372287414
$ws.Cells.Item(24, 10).Value = 8.42220166021302
$ws.Cells.Item(24, 11).Value = 11.18609442946132
$ws.Cells.Item(24, 14).Value = 18.96397037839933

$ws.Cells.Item(25, 2).Value = 11.61145898277239
$ws.Cells.Item(25, 3).Value = 4.893238164320364
$ws.Cells.Item(25, 5).Value = 21.69190660494347
$ws.Cells.Item(25, 6).Value = 43.21761504948566
$ws.Cells.Item(25, 7).Value = 35.528249153894
$ws.Cells.Item(25, 8).Value = 15.89636223146252
$ws.Cells.Item(25, 10).Value = 8.4633168605311
$ws.Cells.Item(25, 11).Value = 10.83310660632164
$ws.Cells.Item(25, 14).Value = 19.08702660574395
